# Commit: "break out stock.yaml completed"
#
# Adds a new "backup" column (R) to the stock-history sheet, backfills it
# with 0 for all existing data rows, resets a couple of stray
# `detect_structure` (Q) flags that had been left over from a prior run,
# flips the `isPivot` (O) flag for the most-recent-at-the-time row, and
# appends two new weekly candles (rows 393-394) that were pulled in after
# the backup column was introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cell R1 = "backup" --------------------------------
# Clone the formatting of the neighbouring header cell (bold, centered,
# bordered) so the new header matches the existing ones, then set text.
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R1").Value = "backup"

# --- 2. Backfill the new "backup" column with 0 for every data row ----
$ws.Range("R2:R392").Value = 0

# --- 3. Row-specific corrections on existing data ----------------------
# detect_structure (Q) got reset to 0 on these two rows.
$ws.Range("Q37").Value = 0
$ws.Range("Q44").Value = 0

# isPivot (O) flips to 1 for row 389.
$ws.Range("O389").Value = 1

# --- 4. Append the two new weekly rows pulled in after the edit --------
# Clone the date-formatted style of column A down into the new rows.
$ws.Range("A392").Copy()
$ws.Range("A393:A394").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("A393").Value = 45460
$ws.Range("B393").Value = 434.8999938964844
$ws.Range("C393").Value = 437.7000122070312
$ws.Range("D393").Value = 426
$ws.Range("E393").Value = 430
$ws.Range("F393").Value = 430
$ws.Range("G393").Value = 4192949
$ws.Range("H393").Value = 2024
$ws.Range("I393").Value = 6
$ws.Range("J393").Value = 17
$ws.Range("K393").Value = 0
$ws.Range("L393").Value = 0
$ws.Range("M393").Value = 0
$ws.Range("N393").Value = 25
$ws.Range("O393").Value = 0
$ws.Range("P393").Value = 0
$ws.Range("Q393").Value = 0

$ws.Range("A394").Value = 45467
$ws.Range("B394").Value = 427
$ws.Range("C394").Value = 433.8999938964844
$ws.Range("D394").Value = 419.5499877929688
$ws.Range("E394").Value = 424.5499877929688
$ws.Range("F394").Value = 424.5499877929688
$ws.Range("G394").Value = 4088819
$ws.Range("H394").Value = 2024
$ws.Range("I394").Value = 6
$ws.Range("J394").Value = 24
$ws.Range("K394").Value = 0
$ws.Range("L394").Value = 0
$ws.Range("M394").Value = 0
$ws.Range("N394").Value = 26
$ws.Range("O394").Value = 0
$ws.Range("P394").Value = 0
$ws.Range("Q394").Value = 0

# The "backup" column has not been computed yet for these brand-new rows,
# so R393/R394 are intentionally left blank (matches the source export's
# empty placeholder for that column).
